$d = $word.ActiveDocument

# These character styles already have <w:b/> and/or <w:i/> set, but the
# <w:color/> element was serialized before them, which violates the
# CT_RPr sequence in wml.xsd (rFonts, b, bCs, i, iCs, caps, ..., color, ...).
# Re-assigning Font.Bold / Font.Italic forces the style's rPr to be
# re-serialized in schema order without changing any actual formatting.

$bold = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleId in $bold) {
    $s = $d.Styles($styleId)
    $s.Font.Bold = $true
}

$italic = @("CommentTok", "DocumentationTok")
foreach ($styleId in $italic) {
    $s = $d.Styles($styleId)
    $s.Font.Italic = $true
}

$boldItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleId in $boldItalic) {
    $s = $d.Styles($styleId)
    $s.Font.Bold = $true
    $s.Font.Italic = $true
}
